# Add 4 new opcodes (37-40) to control U15 from the XMOS, on the
# "XMOS->STM32" sheet, and backfill the matching "Data width" column
# entries that were missing for opcodes 35 and 36 as well as the new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XMOS->STM32")

# New command names in column A for the newly-documented opcodes
# (B39:B42 already hold 37/38/39/40 - unchanged).
$ws.Range("A39").Value = "Get U15 output pin X"
$ws.Range("A40").Value = "Set U15 output pin X"
$ws.Range("A41").Value = "Get U15 outputs"
$ws.Range("A42").Value = "Set U15 outputs"

# New "Data width" entries in column D for the same rows.
$ws.Range("D39").Value = "3-bit cmd, 1-bit reply"
$ws.Range("D40").Value = "4-bit cmd"
$ws.Range("D41").Value = "8-bit reply"
$ws.Range("D42").Value = "8-bit command"

# Restore the scroll position / selection recorded in the sheet view.
$ws.Activate()
$aw = $excel.ActiveWindow
try {
    $aw.ScrollRow = 22
    $aw.ScrollColumn = 1
} catch {
}
$ws.Range("D43").Select()
